# Auto-generated Excel COM-interop script
# Applies the 2024-08-30 daily crime-data update across all affected worksheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 5215
$ws.Range("K3").Value = 5374
$ws.Range("H4").Value = 1741
$ws.Range("K4").Value = 1118
$ws.Range("K5").Value = 387
$ws.Range("K6").Value = 5999
$ws.Range("H7").Value = 26054
$ws.Range("K7").Value = 18093

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K3").Value = 59
$ws.Range("K7").Value = 233

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 335
$ws.Range("K3").Value = 368
$ws.Range("K6").Value = 415
$ws.Range("K7").Value = 1225

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 210
$ws.Range("K3").Value = 286
$ws.Range("K6").Value = 222
$ws.Range("K7").Value = 770

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K3").Value = 104
$ws.Range("K6").Value = 155
$ws.Range("K7").Value = 412

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("K3").Value = 24
$ws.Range("K7").Value = 70

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K7").Value = 537
$ws.Range("K8").Value = 1225
$ws.Range("K11").Value = 347
$ws.Range("K15").Value = 183
$ws.Range("K19").Value = 532
$ws.Range("K20").Value = 417
$ws.Range("K21").Value = 56
$ws.Range("K22").Value = 47
$ws.Range("K23").Value = 189
$ws.Range("K29").Value = 968
$ws.Range("K30").Value = 70
$ws.Range("K31").Value = 196
$ws.Range("K33").Value = 770
$ws.Range("K34").Value = 102
$ws.Range("K36").Value = 238
$ws.Range("K42").Value = 665
$ws.Range("K43").Value = 160
$ws.Range("K47").Value = 120
$ws.Range("K48").Value = 225
$ws.Range("K52").Value = 475
$ws.Range("K53").Value = 233
$ws.Range("K54").Value = 355
$ws.Range("H63").Value = 291
$ws.Range("K63").Value = 52
$ws.Range("K64").Value = 114
$ws.Range("K65").Value = 412
$ws.Range("K67").Value = 690
$ws.Range("K73").Value = 155
$ws.Range("K76").Value = 249
$ws.Range("K77").Value = 128
$ws.Range("K78").Value = 206
$ws.Range("K80").Value = 65
$ws.Range("K84").Value = 134
$ws.Range("K85").Value = 853
$ws.Range("K89").Value = 263
$ws.Range("K91").Value = 195
$ws.Range("K94").Value = 238
$ws.Range("K96").Value = 194
$ws.Range("K97").Value = 145
$ws.Range("H101").Value = 26054
$ws.Range("K101").Value = 18093

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K3").Value = 46
$ws.Range("K7").Value = 196

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K3").Value = 245
$ws.Range("K6").Value = 194
$ws.Range("K7").Value = 690

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K2").Value = 43
$ws.Range("K7").Value = 134

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K6").Value = 189
$ws.Range("K7").Value = 355

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 277
$ws.Range("K3").Value = 348
$ws.Range("K6").Value = 268
$ws.Range("K7").Value = 968

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K3").Value = 52
$ws.Range("K7").Value = 225

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K6").Value = 170
$ws.Range("K7").Value = 532

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K2").Value = 53
$ws.Range("K5").Value = 2
$ws.Range("K7").Value = 249

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 176
$ws.Range("K3").Value = 205
$ws.Range("K4").Value = 26
$ws.Range("K5").Value = 7
$ws.Range("K7").Value = 665

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K3").Value = 46
$ws.Range("K6").Value = 74
$ws.Range("K7").Value = 206

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K2").Value = 52
$ws.Range("K3").Value = 68
$ws.Range("K6").Value = 51
$ws.Range("K7").Value = 189

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K6").Value = 84
$ws.Range("K7").Value = 194

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K2").Value = 49
$ws.Range("K7").Value = 195

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("K5").Value = 2
$ws.Range("K7").Value = 56

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("K3").Value = 34
$ws.Range("K7").Value = 114

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K3").Value = 135
$ws.Range("K7").Value = 417

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K3").Value = 69
$ws.Range("K7").Value = 238

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K6").Value = 142
$ws.Range("K7").Value = 537

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("K2").Value = 36
$ws.Range("K7").Value = 102

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K3").Value = 45
$ws.Range("K6").Value = 102
$ws.Range("K7").Value = 238

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("K3").Value = 35
$ws.Range("K6").Value = 40
$ws.Range("K7").Value = 120

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("K3").Value = 44
$ws.Range("K6").Value = 58
$ws.Range("K7").Value = 183

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K2").Value = 117
$ws.Range("K7").Value = 347

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K2").Value = 49
$ws.Range("K7").Value = 155

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("K3").Value = 27
$ws.Range("K6").Value = 85
$ws.Range("K7").Value = 145

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K6").Value = 79
$ws.Range("K7").Value = 263

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K4").Value = 21
$ws.Range("K6").Value = 64
$ws.Range("K7").Value = 160

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 285
$ws.Range("K3").Value = 285
$ws.Range("K6").Value = 207
$ws.Range("K7").Value = 853

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("K3").Value = 15
$ws.Range("K6").Value = 7
$ws.Range("K7").Value = 47

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("K2").Value = 56
$ws.Range("K7").Value = 128

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("K6").Value = 30
$ws.Range("K7").Value = 65

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K3").Value = 133
$ws.Range("K5").Value = 17
$ws.Range("K7").Value = 475
